$wb = $excel.ActiveWorkbook

# ALC!row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 109.333336
$ws.Range("I39").Value = 71.2
$ws.Range("K39").Value = 213.6
$ws.Range("M39").Value = 82.39999999999998

# ALC!row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 73165.10000000001
$ws.Range("I62").Value = 42781
$ws.Range("J62").Value = 103549.2
$ws.Range("K62").Value = 42781
$ws.Range("L62").Value = 103549.2
$ws.Range("M62").Value = -42157
$ws.Range("N62").Value = -104797.2

# ALC!row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 1319047.2
$ws.Range("I64").Value = 1757209.1
$ws.Range("J64").Value = 4561.5
$ws.Range("K64").Value = 1757209.1
$ws.Range("L64").Value = 4561.5
$ws.Range("M64").Value = -1756961.1
$ws.Range("N64").Value = -5057.5

# ALC!row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 73165.10000000001
$ws.Range("I65").Value = 42781
$ws.Range("J65").Value = 103549.2
$ws.Range("K65").Value = 213905
$ws.Range("L65").Value = 517746
$ws.Range("M65").Value = -210785
$ws.Range("N65").Value = -523986

# ALC!row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 1319047.2
$ws.Range("I67").Value = 1757209.1
$ws.Range("J67").Value = 4561.5
$ws.Range("K67").Value = 1757209.1
$ws.Range("L67").Value = 4561.5
$ws.Range("M67").Value = -1756351.1
$ws.Range("N67").Value = -6277.5

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1323.6666
$ws.Range("I137").Value = 811.2
$ws.Range("J137").Value = 1579.9
$ws.Range("K137").Value = 2433.6
$ws.Range("L137").Value = 4739.700000000001
$ws.Range("M137").Value = 116.3999999999996
$ws.Range("N137").Value = -9839.700000000001

# ARM!row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 500
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -732

# ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4070.09
$ws.Range("I32").Value = 3195.593
$ws.Range("J32").Value = 9442
$ws.Range("K32").Value = 3195.593
$ws.Range("L32").Value = 9442
$ws.Range("M32").Value = -2908.593
$ws.Range("N32").Value = -10016

# CRP!row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56.75
$ws.Range("I7").Value = 53.444443
$ws.Range("J7").Value = 66.666664
$ws.Range("K7").Value = 53.444443
$ws.Range("L7").Value = 66.666664
$ws.Range("M7").Value = 59.555557
$ws.Range("N7").Value = -292.666664

# CRP!row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 898.0769
$ws.Range("I22").Value = 925.4286
$ws.Range("J22").Value = 866.1667
$ws.Range("K22").Value = 925.4286
$ws.Range("L22").Value = 866.1667
$ws.Range("M22").Value = -575.4286
$ws.Range("N22").Value = -1566.1667

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3436.4138
$ws.Range("I31").Value = 2853
$ws.Range("J31").Value = 3980.9333
$ws.Range("K31").Value = 2853
$ws.Range("L31").Value = 3980.9333
$ws.Range("M31").Value = -2558
$ws.Range("N31").Value = -4570.933300000001

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3436.4138
$ws.Range("I34").Value = 2853
$ws.Range("J34").Value = 3980.9333
$ws.Range("K34").Value = 2853
$ws.Range("L34").Value = 3980.9333
$ws.Range("M34").Value = -2651
$ws.Range("N34").Value = -4384.933300000001

# CRP!row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17739.8
$ws.Range("J68").Value = 17739.8
$ws.Range("L68").Value = 17739.8
$ws.Range("N68").Value = -19237.8

# CRP!row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17739.8
$ws.Range("J71").Value = 17739.8
$ws.Range("L71").Value = 53219.39999999999
$ws.Range("N71").Value = -60707.39999999999

# CUL!row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 147
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 171.25
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 513.75
$ws.Range("M7").Value = -38
$ws.Range("N7").Value = -737.75

# CUL!row 16
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1371.6666
$ws.Range("J16").Value = 5980
$ws.Range("L16").Value = 17940
$ws.Range("N16").Value = -18286

# CUL!row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2060.5881
$ws.Range("J80").Value = 2076.875
$ws.Range("L80").Value = 6230.625
$ws.Range("N80").Value = -8102.625

# CUL!row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 2060.5881
$ws.Range("J83").Value = 2076.875
$ws.Range("L83").Value = 18691.875
$ws.Range("N83").Value = -28051.875

# GSM!row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3562
$ws.Range("I80").Value = 2888.3333
$ws.Range("J80").Value = 4011.111
$ws.Range("K80").Value = 2888.3333
$ws.Range("L80").Value = 4011.111
$ws.Range("M80").Value = -1890.3333
$ws.Range("N80").Value = -6007.111

# GSM!row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3562
$ws.Range("I83").Value = 2888.3333
$ws.Range("J83").Value = 4011.111
$ws.Range("K83").Value = 14441.6665
$ws.Range("L83").Value = 20055.555
$ws.Range("M83").Value = -9449.666499999999
$ws.Range("N83").Value = -30039.555

# LTW!row 62
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# LTW!row 65
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# LTW!row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2184.1924
$ws.Range("J68").Value = 2162.6365
$ws.Range("L68").Value = 2162.6365
$ws.Range("N68").Value = -3660.6365

# LTW!row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2184.1924
$ws.Range("J71").Value = 2162.6365
$ws.Range("L71").Value = 10813.1825
$ws.Range("N71").Value = -18301.1825

# WVR!row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3272.8333
$ws.Range("I62").Value = 2567
$ws.Range("J62").Value = 3625.75
$ws.Range("K62").Value = 2567
$ws.Range("L62").Value = 3625.75
$ws.Range("M62").Value = -1943
$ws.Range("N62").Value = -4873.75

# WVR!row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3272.8333
$ws.Range("I65").Value = 2567
$ws.Range("J65").Value = 3625.75
$ws.Range("K65").Value = 12835
$ws.Range("L65").Value = 18128.75
$ws.Range("M65").Value = -9715
$ws.Range("N65").Value = -24368.75

# WVR!row 68
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101622

# WVR!row 71
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -308112

# WVR!row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1589.2222
$ws.Range("I122").Value = 1471.8572
$ws.Range("K122").Value = 4415.571599999999
$ws.Range("M122").Value = -1965.571599999999
